$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text (inline/shared string) even when
# it looks numeric or date-like, without introducing any new cell style.
# We do this by writing a text formula that evaluates to the exact string,
# then converting the cell to a static value via copy / paste-special values.
function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $escaped = $Text -replace '"', '""'
    $ws.Range($Address).Formula = "=""" + $escaped + """"
    $ws.Range($Address).Copy()
    $ws.Range($Address).PasteSpecial(-4163)
}

# --- Fix row 18: B18 must be a text cell containing "1517492" (was numeric) ---
Set-TextValue "B18" "1517492"

# --- Add new row 19 ---
$ws.Range("A19").Value = 123444
$ws.Range("B19").Value = "PROPRO"
$ws.Range("C19").Value = "uuuuuuuuuuu"
$ws.Range("D19").Value = "Mètre"
$ws.Range("E19").Value = "Barre de 6m"
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = "Site principal"
$ws.Range("I19").Value = "Stockage"
$ws.Range("J19").Value = "E3"
$ws.Range("K19").Value = "FournX"
$ws.Range("L19").Value = 45
$ws.Range("M19").Value = "Profilés"
$ws.Range("N19").Value = "Structure"
Set-TextValue "O19" "2481023879"
$ws.Range("P19").Value = 40
Set-TextValue "Q19" "2025-05-28"

$wb.Save()
